$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# New Cypher query text for the Participants tab (column B, row 2)
$newParticipantQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE g.platform in ['Illumina Next Seq 550']
with p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id LIMIT 100
'@

$ws.Range("B2").Value = $newParticipantQuery

# Recalculate the row height to fit the new (longer) wrapped text
$ws.Rows.Item(2).RowHeight = 279

# Update selection to B5 (matches target sheetView selection)
$ws.Range("B5").Select()
